# Edit the opening paragraph of the ΕΚΘΕΣΗ ΕΞΕΤΑΣΗΣ ΚΑΤΗΓΟΡΟΥΜΕΝΟΥ document:
# the template's opening sentence had its "Στην ___" gap filled with a town
# name ("Θέρμη") followed by a tab, and the rest of the blank template
# gaps filled in with the case's details.
$d = $word.ActiveDocument

# --- Edit 1 --------------------------------------------------------------
# Locate the opening paragraph by its (old) literal text, then replace it
# with three runs: "Στην Θέρμη", a real tab character (<w:tab/>), and the
# remainder of the sentence with all of its blanks filled in. InsertXML is
# used (after deleting the located text) so the tab becomes a genuine
# <w:tab/> run, matching what Word records for an actually-typed Tab key.
$target = $d.Content
$found1 = $target.Find.Execute('Στην  σήμερα την  του μήνα  του έτους  ημέρα εβδομάδας  και ώρα  ενώπιον εμού του  του  Θεσσαλονίκης, παρισταμένου  και του   της ιδίας υπηρεσίας, που προσλήφθηκε ως Β'' Ανακριτικός Υπάλληλος, προσκλήθηκε ο  κατωτέρω σημειούμενος  κατηγορούμενος, που ονομάζεται    του  και της  γεν.  στη  κατ.,αριθμός τηλεφώνου , ηλεκτρονικό  ταχυδρομείου, κάτοχος του υπ αριθμόν  που εκδόθηκε την  από Α.Φ.Μ : , Δ.Ο.Υ : ,στον οποίο γνωστοποιήσαμε ότι κατηγορείται για παράβαση του/των άρθρων  [] τπυ Π.Κ. και εξηγήσαμε με σαφήνεια και πληρότητα σ’ αυτόν βάσει του άρθρου 95 του Κώδικα Ποινικής  Δικονομίας όλα τα εκ των  άρθρων 91,95,96,97,98,99,100,103 και 104  του  Κ.Π.Δ.  δικαιώματά  του/της  και  αναλυτικότερα :', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the opening paragraph text to replace"
}
$target.Delete()
$xmlPayload1 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Στην Θέρμη</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve"> σήμερα την ψσδψδσ του μήνα σδψδσ του έτους ψδσψδσψ ημέρα εβδομάδας ψσδψδσψδσψδσψ και ώρα 1029 ενώπιον εμού του Υ/Α Αθαθααθδσψνβσδκνψ του ψσδκψωδφσωνδφω Θεσσαλονίκης, παρισταμένου  και του  ωσωσδωσδω της ιδίας υπηρεσίας, που προσλήφθηκε ως Β'' Ανακριτικός Υπάλληλος, προσκλήθηκε ο  κατωτέρω σημειούμενος  κατηγορούμενος, που ονομάζεται  ωδφωω ωδφωδδφω του ωδφωδφωδφ και της ωδωδφωφδ γεν. ωδωδφωφδδφω στη δφωφωφδ κατ.ωδφωφδωφδ,αριθμός τηλεφώνου ωφδφωφδω, ηλεκτρονικό  ταχυδρομείουωδφωδφωδφ, κάτοχος του υπ αριθμόν ωδωδφ που εκδόθηκε την ωδωδφωδφωφδ από δφσφσδσφδσΑ.Φ.Μ : φσδφδσφδσφδσ, Δ.Ο.Υ : φσδφδσφδσφ,στον οποίο γνωστοποιήσαμε ότι κατηγορείται για παράβαση του/των άρθρων  [''308 "Σωματικεσ"'', ''361 ερρωηφ''] τπυ Π.Κ. και εξηγήσαμε με σαφήνεια και πληρότητα σ’ αυτόν βάσει του άρθρου 95 του Κώδικα Ποινικής  Δικονομίας όλα τα εκ των  άρθρων 91,95,96,97,98,99,100,103 και 104  του  Κ.Π.Δ.  δικαιώματά  του/της  και  αναλυτικότερα :</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xmlPayload1)

# --- Edit 2 --------------------------------------------------------------
# Fill in the two time blanks in the closing "Η παρούσα έκθεση ..." line.
# Re-locate via Find, then set Paragraph.Range.Text directly (rather than
# Find.Execute's Replace) so xml:space="preserve" is kept on the run even
# though the inserted value only leaves an internal (non-edge) double
# space, matching the original run's whitespace-preserving markup.
$target2 = $d.Content
$found2 = $target2.Find.Execute('Η παρούσα έκθεση άρχισε να συντάσσεται την  ώρα και περατώθηκε την   ώρα', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the closing paragraph text to replace"
}
$para2 = $target2.Paragraphs(1)
$para2.Range.Text = 'Η παρούσα έκθεση άρχισε να συντάσσεται την 1029 ώρα και περατώθηκε την 32423  ώρα'

Write-Output "Edits applied"
